$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the grid with two new columns (H, I) and widen columns G/H to match B-F ---
# ColumnWidth is quantized to pixel steps by the engine, so this is the closest achievable
# value to the original custom width (10.7109375) used by columns B-F.
$ws.Columns.Item(7).ColumnWidth = 9.83
$ws.Columns.Item(8).ColumnWidth = 9.83

# Force creation of the new (but blank) header cells G1, H1, G2, H2 with the same right-aligned
# style ("s=4") already used by the rest of row 1 / row 2, mirroring existing cells like F1/F2.
$ws.Range("G1").HorizontalAlignment = -4152
$ws.Range("H1").HorizontalAlignment = -4152
$ws.Range("G2").HorizontalAlignment = -4152
$ws.Range("H2").HorizontalAlignment = -4152

# --- Write the full A1:I27 grid as text (values are stored as literal strings in the sheet, e.g.
# dates like "28-07-2024" and comma-decimal numbers like "495,36") ---
$ws.Range("A1").Value = "Kurs (PLN/MWh)"
$ws.Range("B1").Value = ""
$ws.Range("C1").Value = ""
$ws.Range("D1").Value = ""
$ws.Range("E1").Value = ""
$ws.Range("F1").Value = ""
$ws.Range("G1").Value = ""
$ws.Range("H1").Value = ""
$ws.Range("A2").Value = "FIXING I"
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("A3").Value = ""
$ws.Range("B3").Value = "28-07-2024"
$ws.Range("C3").Value = "29-07-2024"
$ws.Range("D3").Value = "30-07-2024"
$ws.Range("E3").Value = "31-07-2024"
$ws.Range("F3").Value = "01-08-2024"
$ws.Range("G3").Value = "02-08-2024"
$ws.Range("H3").Value = "03-08-2024"
$ws.Range("I3").Value = "04-08-2024"
$ws.Range("A4").Value = "0-1"
$ws.Range("B4").Value = "495,36"
$ws.Range("C4").Value = "317,99"
$ws.Range("D4").Value = "433,00"
$ws.Range("E4").Value = "473,14"
$ws.Range("F4").Value = "450,00"
$ws.Range("G4").Value = "465,00"
$ws.Range("H4").Value = "490,00"
$ws.Range("I4").Value = "442,00"
$ws.Range("A5").Value = "1-2"
$ws.Range("B5").Value = "439,80"
$ws.Range("C5").Value = "280,00"
$ws.Range("D5").Value = "380,00"
$ws.Range("E5").Value = "418,00"
$ws.Range("F5").Value = "400,00"
$ws.Range("G5").Value = "400,00"
$ws.Range("H5").Value = "450,00"
$ws.Range("I5").Value = "412,45"
$ws.Range("A6").Value = "2-3"
$ws.Range("B6").Value = "399,00"
$ws.Range("C6").Value = "251,00"
$ws.Range("D6").Value = "389,50"
$ws.Range("E6").Value = "406,26"
$ws.Range("F6").Value = "390,00"
$ws.Range("G6").Value = "392,00"
$ws.Range("H6").Value = "424,41"
$ws.Range("I6").Value = "384,00"
$ws.Range("A7").Value = "3-4"
$ws.Range("B7").Value = "400,00"
$ws.Range("C7").Value = "250,00"
$ws.Range("D7").Value = "376,00"
$ws.Range("E7").Value = "407,00"
$ws.Range("F7").Value = "390,00"
$ws.Range("G7").Value = "388,00"
$ws.Range("H7").Value = "436,90"
$ws.Range("I7").Value = "375,00"
$ws.Range("A8").Value = "4-5"
$ws.Range("B8").Value = "394,00"
$ws.Range("C8").Value = "272,00"
$ws.Range("D8").Value = "381,08"
$ws.Range("E8").Value = "408,00"
$ws.Range("F8").Value = "397,20"
$ws.Range("G8").Value = "390,00"
$ws.Range("H8").Value = "428,00"
$ws.Range("I8").Value = "386,00"
$ws.Range("A9").Value = "5-6"
$ws.Range("B9").Value = "382,00"
$ws.Range("C9").Value = "298,00"
$ws.Range("D9").Value = "393,55"
$ws.Range("E9").Value = "443,55"
$ws.Range("F9").Value = "437,00"
$ws.Range("G9").Value = "431,28"
$ws.Range("H9").Value = "447,91"
$ws.Range("I9").Value = "392,00"
$ws.Range("A10").Value = "6-7"
$ws.Range("B10").Value = "412,00"
$ws.Range("C10").Value = "406,00"
$ws.Range("D10").Value = "509,11"
$ws.Range("E10").Value = "550,00"
$ws.Range("F10").Value = "530,00"
$ws.Range("G10").Value = "542,00"
$ws.Range("H10").Value = "480,00"
$ws.Range("I10").Value = "396,00"
$ws.Range("A11").Value = "7-8"
$ws.Range("B11").Value = "390,00"
$ws.Range("C11").Value = "378,00"
$ws.Range("D11").Value = "506,40"
$ws.Range("E11").Value = "530,00"
$ws.Range("F11").Value = "520,00"
$ws.Range("G11").Value = "540,00"
$ws.Range("H11").Value = "488,68"
$ws.Range("I11").Value = "359,00"
$ws.Range("A12").Value = "8-9"
$ws.Range("B12").Value = "361,04"
$ws.Range("C12").Value = "299,00"
$ws.Range("D12").Value = "416,30"
$ws.Range("E12").Value = "420,00"
$ws.Range("F12").Value = "432,41"
$ws.Range("G12").Value = "488,50"
$ws.Range("H12").Value = "524,27"
$ws.Range("I12").Value = "259,99"
$ws.Range("A13").Value = "9-10"
$ws.Range("B13").Value = "269,99"
$ws.Range("C13").Value = "137,99"
$ws.Range("D13").Value = "283,99"
$ws.Range("E13").Value = "303,99"
$ws.Range("F13").Value = "355,57"
$ws.Range("G13").Value = "450,00"
$ws.Range("H13").Value = "406,00"
$ws.Range("I13").Value = "101,00"
$ws.Range("A14").Value = "10-11"
$ws.Range("B14").Value = "100,00"
$ws.Range("C14").Value = "0,00"
$ws.Range("D14").Value = "180,11"
$ws.Range("E14").Value = "198,12"
$ws.Range("F14").Value = "288,00"
$ws.Range("G14").Value = "404,00"
$ws.Range("H14").Value = "388,80"
$ws.Range("I14").Value = "30,00"
$ws.Range("A15").Value = "11-12"
$ws.Range("B15").Value = "0,00"
$ws.Range("C15").Value = "-31,00"
$ws.Range("D15").Value = "90,11"
$ws.Range("E15").Value = "105,95"
$ws.Range("F15").Value = "263,40"
$ws.Range("G15").Value = "390,34"
$ws.Range("H15").Value = "377,76"
$ws.Range("I15").Value = "45,99"
$ws.Range("A16").Value = "12-13"
$ws.Range("B16").Value = "-30,01"
$ws.Range("C16").Value = "-55,00"
$ws.Range("D16").Value = "30,00"
$ws.Range("E16").Value = "70,00"
$ws.Range("F16").Value = "235,90"
$ws.Range("G16").Value = "387,59"
$ws.Range("H16").Value = "361,40"
$ws.Range("I16").Value = "79,99"
$ws.Range("A17").Value = "13-14"
$ws.Range("B17").Value = "-71,01"
$ws.Range("C17").Value = "-54,01"
$ws.Range("D17").Value = "10,11"
$ws.Range("E17").Value = "70,00"
$ws.Range("F17").Value = "241,99"
$ws.Range("G17").Value = "394,00"
$ws.Range("H17").Value = "326,70"
$ws.Range("I17").Value = "199,89"
$ws.Range("A18").Value = "14-15"
$ws.Range("B18").Value = "-190,00"
$ws.Range("C18").Value = "-30,01"
$ws.Range("D18").Value = "21,40"
$ws.Range("E18").Value = "100,00"
$ws.Range("F18").Value = "233,99"
$ws.Range("G18").Value = "406,00"
$ws.Range("H18").Value = "300,00"
$ws.Range("I18").Value = "217,99"
$ws.Range("A19").Value = "15-16"
$ws.Range("B19").Value = "-176,00"
$ws.Range("C19").Value = "0,01"
$ws.Range("D19").Value = "165,00"
$ws.Range("E19").Value = "239,99"
$ws.Range("F19").Value = "295,67"
$ws.Range("G19").Value = "482,86"
$ws.Range("H19").Value = "340,00"
$ws.Range("I19").Value = "330,27"
$ws.Range("A20").Value = "16-17"
$ws.Range("B20").Value = "-111,00"
$ws.Range("C20").Value = "90,00"
$ws.Range("D20").Value = "315,99"
$ws.Range("E20").Value = "370,60"
$ws.Range("F20").Value = "385,69"
$ws.Range("G20").Value = "482,86"
$ws.Range("H20").Value = "399,99"
$ws.Range("I20").Value = "386,00"
$ws.Range("A21").Value = "17-18"
$ws.Range("B21").Value = "-15,00"
$ws.Range("C21").Value = "281,99"
$ws.Range("D21").Value = "380,60"
$ws.Range("E21").Value = "379,00"
$ws.Range("F21").Value = "444,97"
$ws.Range("G21").Value = "469,00"
$ws.Range("H21").Value = "492,90"
$ws.Range("I21").Value = "482,00"
$ws.Range("A22").Value = "18-19"
$ws.Range("B22").Value = "159,30"
$ws.Range("C22").Value = "380,00"
$ws.Range("D22").Value = "481,00"
$ws.Range("E22").Value = "499,40"
$ws.Range("F22").Value = "503,60"
$ws.Range("G22").Value = "549,00"
$ws.Range("H22").Value = "490,00"
$ws.Range("I22").Value = "465,00"
$ws.Range("A23").Value = "19-20"
$ws.Range("B23").Value = "343,40"
$ws.Range("C23").Value = "558,00"
$ws.Range("D23").Value = "682,00"
$ws.Range("E23").Value = "632,00"
$ws.Range("F23").Value = "612,00"
$ws.Range("G23").Value = "600,00"
$ws.Range("H23").Value = "560,00"
$ws.Range("I23").Value = "489,00"
$ws.Range("A24").Value = "20-21"
$ws.Range("B24").Value = "420,80"
$ws.Range("C24").Value = "680,00"
$ws.Range("D24").Value = "853,40"
$ws.Range("E24").Value = "825,00"
$ws.Range("F24").Value = "738,99"
$ws.Range("G24").Value = "670,00"
$ws.Range("H24").Value = "600,00"
$ws.Range("I24").Value = "537,44"
$ws.Range("A25").Value = "21-22"
$ws.Range("B25").Value = "445,80"
$ws.Range("C25").Value = "625,00"
$ws.Range("D25").Value = "720,00"
$ws.Range("E25").Value = "634,00"
$ws.Range("F25").Value = "585,00"
$ws.Range("G25").Value = "553,90"
$ws.Range("H25").Value = "586,61"
$ws.Range("I25").Value = "512,20"
$ws.Range("A26").Value = "22-23"
$ws.Range("B26").Value = "390,00"
$ws.Range("C26").Value = "496,00"
$ws.Range("D26").Value = "549,00"
$ws.Range("E26").Value = "520,00"
$ws.Range("F26").Value = "537,00"
$ws.Range("G26").Value = "513,70"
$ws.Range("H26").Value = "514,80"
$ws.Range("I26").Value = "473,61"
$ws.Range("A27").Value = "23-24"
$ws.Range("B27").Value = "338,80"
$ws.Range("C27").Value = "428,90"
$ws.Range("D27").Value = "458,00"
$ws.Range("E27").Value = "443,80"
$ws.Range("F27").Value = "450,00"
$ws.Range("G27").Value = "449,00"
$ws.Range("H27").Value = "441,87"
$ws.Range("I27").Value = "437,39"
